$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# PHASE 0: detach the existing comment from Q1 (it will be re-created at AA1,
# the new home of Q1's content, once the old cells have been relocated).
# ---------------------------------------------------------------------------
$ws.Range("Q1").Comment.Delete()

# ---------------------------------------------------------------------------
# PHASE 1: relocate the pre-existing trailing cells (date + string samples)
# out of the way of the new CURRENCY columns being inserted at O1:Y1.
# Order matters: R1 is both a source and (later) a destination, so read it
# out first.
# ---------------------------------------------------------------------------
$ws.Range("R1").Cut($ws.Range("AB1"))      # old R1 (string, red font)  -> AB1
$ws.Range("O1").Cut($ws.Range("Z1"))       # old O1 (date)              -> Z1
$ws.Range("Q1").Cut($ws.Range("AA1"))      # old Q1 (string)            -> AA1
$ws.Range("S1").Cut($ws.Range("AC1"))      # old S1 (string, Automatic) -> AC1
$ws.Range("P1").Cut($ws.Range("R1"))       # old P1 (currency, red)     -> R1

# Recreate the comment at its new home.
$ws.Range("AA1").AddComment("Note")

# ---------------------------------------------------------------------------
# PHASE 2: populate the new CURRENCY ('R$') demonstration cells, O1:Y1.
# ---------------------------------------------------------------------------
$ws.Range("O1").NumberFormat = """R$""\ #,##0;\-""R$""\ #,##0"
$ws.Range("O1").Value = 1

$ws.Range("P1").NumberFormat = """R$""\ #,##0;[Red]\-""R$""\ #,##0"
$ws.Range("P1").Value = 1

$ws.Range("Q1").NumberFormat = """R$""\ #,##0.00;[Red]""R$""\ #,##0.00"
$ws.Range("Q1").Value = 1

# R1 already holds the relocated old-P1 value (1) with its original
# "R$ #,##0.00;[Red]-R$ #,##0.00" currency format - nothing further to do.

$ws.Range("S1").NumberFormat = "_-""R$""\ * #,##0_-;\-""R$""\ * #,##0_-;_-""R$""\ * ""-""_-;_-@_-"
$ws.Range("S1").Value = 1

$ws.Range("T1").NumberFormat = "_-""R$""\ * #,##0.00_-;\-""R$""\ * #,##0.00_-;_-""R$""\ * ""-""??_-;_-@_-"
$ws.Range("T1").Value = 1

$ws.Range("U1").NumberFormat = """R$""\ #,##0.00;[Red]""R$""\ #,##0.00"
$ws.Range("U1").Formula = "=T1+T1"

$ws.Range("V1").NumberFormat = "@"
$ws.Range("V1").Value = "R$ 1.00"

$ws.Range("W1").Formula = "=CONCAT(""R$ "",""1.00"")"

$ws.Range("X1").NumberFormat = "@"
$ws.Range("X1").Value = "R$ 1,00"

$ws.Range("Y1").Formula = "=CONCAT(""R$ "",""1,00"")"
